$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("private")

# Sort the A2:C96 table by Type (column C) then Owner (column A), ascending,
# matching the workbook's "Sort" feature (adds a persisted sortState too).
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("C2:C96"))
$ws.Sort.SortFields.Add($ws.Range("A2:A96"))
$ws.Sort.SetRange($ws.Range("A1:C96"))
$ws.Sort.Header = 1
$ws.Sort.Apply()

# Widen the columns to fit the (now sorted) content.
$ws.Columns.Item(1).ColumnWidth = 37.5
$ws.Columns.Item(2).ColumnWidth = 34.33333333333333
$ws.Columns.Item(3).ColumnWidth = 70.66666666666667

# Make "private" the active sheet/tab, with C6 as the selected cell.
$ws.Activate()
$ws.Range("C6").Select()
